$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1 previously held 0.0; update it to 10.0 (matching C1's value)
$ws.Range("D1").Value = 10.0
